# Update stat values on Sheet1 to reflect the refreshed tennis stats
# (cfs_6_0.3.xlsx data refresh, "Add files via upload").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 0.501
$ws.Range("G3").Value = 0.8080000000000001
$ws.Range("F4").Value = 0.58
$ws.Range("G4").Value = 0.734
$ws.Range("F5").Value = 0.627
$ws.Range("L5").Value = 0.429
$ws.Range("H6").Value = 0.404
$ws.Range("J8").Value = 0.431
$ws.Range("G9").Value = 0.744
$ws.Range("K10").Value = 0.418
$ws.Range("J13").Value = 0.384
$ws.Range("K13").Value = 0.458
$ws.Range("K16").Value = 0.459
$ws.Range("E32").Value = 0.524
$ws.Range("G49").Value = 0.599
$ws.Range("G62").Value = 0.616
$ws.Range("I62").Value = 0.537
$ws.Range("G65").Value = 0.74
$ws.Range("L65").Value = 0.447
$ws.Range("G68").Value = 0.624
$ws.Range("J68").Value = 0.344
$ws.Range("E69").Value = 0.542
$ws.Range("F69").Value = 0.494
$ws.Range("E77").Value = 0.601
$ws.Range("G77").Value = 0.711
$ws.Range("I81").Value = 0.58
$ws.Range("J81").Value = 0.418
$ws.Range("I84").Value = 0.58
$ws.Range("J84").Value = 0.418
$ws.Range("G88").Value = 0.671
$ws.Range("E89").Value = 0.52
$ws.Range("J89").Value = 0.259
$ws.Range("F90").Value = 0.477
$ws.Range("K90").Value = 0.383
$ws.Range("F91").Value = 0.509
$ws.Range("H91").Value = 0.362
$ws.Range("L93").Value = 0.395
$ws.Range("K97").Value = 0.337
$ws.Range("L97").Value = 0.387
$ws.Range("E99").Value = 0.523
$ws.Range("I99").Value = 0.501
$ws.Range("F100").Value = 0.477
$ws.Range("K100").Value = 0.383
